# Weekly update: a new "Coliflor" price record for
# Vega Central Mapocho de Santiago is inserted as row 489 (pushing the
# existing rows 489:520 down to 490:521). All other rows are left intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at 489, shifting rows 489:520 -> 490:521.
$ws.Rows.Item(489).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A489").Value = 9
$ws.Range("B489").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C489").Value = "Metropolitana"
$ws.Range("D489").Value = 44585
$ws.Range("E489").Value = 13
$ws.Range("F489").Value = 100112008
$ws.Range("G489").Value = "Coliflor"
$ws.Range("H489").Value = "Sin especificar"
$ws.Range("I489").Value = "Primera"
$ws.Range("J489").Value = 1600
$ws.Range("K489").Value = 900
$ws.Range("L489").Value = 950
$ws.Range("M489").Value = 925
$ws.Range("N489").Value = "`$/unidad"
$ws.Range("O489").Value = "Región Metropolitana"
$ws.Range("P489").Value = 925
$ws.Range("Q489").Value = 1
$ws.Range("R489").Value = "Hortaliza"
